$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (pushes existing rows 8-55 down to 9-56)
$ws.Rows("8:8").Insert()

# Populate the newly inserted row with the new localization entry
$ws.Range("A8").Value = "cycleEnd"
$ws.Range("B8").Value = "CYCLE COMPLETE"

# Update selection / view to match the authored state (selection on B8, no frozen topLeftCell scroll)
$ws.Range("B8").Select()
